$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 870, pushing all
# existing rows (870..908) down to (872..910). This mirrors the weekly
# update that prepends a new price observation block to the data set.
$ws.Rows.Item(870).Insert()
$ws.Rows.Item(871).Insert()

# Populate the two newly-inserted rows with the new weekly record
# (Lechuga / Escarola, Primera & Segunda) for date 44753 (2022-07-11).

# Row 870: Escarola / Primera
$ws.Cells.Item(870, 1).Value = 1
$ws.Cells.Item(870, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(870, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(870, 4).Value = 44753
$ws.Cells.Item(870, 5).Value = 15
$ws.Cells.Item(870, 6).Value = 100112033
$ws.Cells.Item(870, 7).Value = "Lechuga"
$ws.Cells.Item(870, 8).Value = "Escarola"
$ws.Cells.Item(870, 9).Value = "Primera"
$ws.Cells.Item(870, 10).Value = 120
$ws.Cells.Item(870, 11).Value = 9000
$ws.Cells.Item(870, 12).Value = 10000
$ws.Cells.Item(870, 13).Value = 9500
$ws.Cells.Item(870, 14).Value = "$/caja 12 unidades"
$ws.Cells.Item(870, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(870, 16).Value = 792
$ws.Cells.Item(870, 17).Value = 12
$ws.Cells.Item(870, 18).Value = "Hortaliza"

# Row 871: Escarola / Segunda
$ws.Cells.Item(871, 1).Value = 1
$ws.Cells.Item(871, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(871, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(871, 4).Value = 44753
$ws.Cells.Item(871, 5).Value = 15
$ws.Cells.Item(871, 6).Value = 100112033
$ws.Cells.Item(871, 7).Value = "Lechuga"
$ws.Cells.Item(871, 8).Value = "Escarola"
$ws.Cells.Item(871, 9).Value = "Segunda"
$ws.Cells.Item(871, 10).Value = 120
$ws.Cells.Item(871, 11).Value = 9000
$ws.Cells.Item(871, 12).Value = 10000
$ws.Cells.Item(871, 13).Value = 9500
$ws.Cells.Item(871, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(871, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(871, 16).Value = 528
$ws.Cells.Item(871, 17).Value = 18
$ws.Cells.Item(871, 18).Value = "Hortaliza"
